# Updated cryptos list on Tue Sep 12 15:41:22 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.280.88"
$ws.Range("E2").Value = "  +3.81%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.607.72"
$ws.Range("E3").Value = "  +2.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.62%  "

# Row 5 - BNB (price is a plain number-looking string -> keep as text)
$ws.Range("D5").Value = "'212.96"
$ws.Range("E5").Value = "  +2.41%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.65%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.486"
$ws.Range("E7").Value = "  +1.75%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +2.22%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.0619"
$ws.Range("E9").Value = "  +1.79%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'18.12"
$ws.Range("E10").Value = "  +0.75%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0817"
$ws.Range("E11").Value = "  +4.45%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.830.86"
$ws.Range("E12").Value = "  +2.22%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.603.41"
$ws.Range("E13").Value = "  +1.96%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.42%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.30%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.241.44"
$ws.Range("E16").Value = "  +3.61%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'60.73"
$ws.Range("E17").Value = "  +1.67%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +2.27%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.54%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'198.79"
$ws.Range("E20").Value = "  +7.11%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +2.62%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "'9.37"
$ws.Range("E22").Value = "  +0.35%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +1.86%  "

# Row 24 - Monero
$ws.Range("D24").Value = "'142.69"
$ws.Range("E24").Value = "  +1.27%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +3.26%  "

# Rows 26 & 27 swap places: Stellar <-> BinanceUSD
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.126"
$ws.Range("E27").Value = "  -1.85%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'15.18"
$ws.Range("E28").Value = "  +2.03%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "'6.48"
$ws.Range("E29").Value = "  +0.11%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.48%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +2.37%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.80%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.11%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.96%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'2.35"
$ws.Range("E35").Value = "  +4.96%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.107.61"
$ws.Range("E36").Value = "  +1.85%  "

# Row 37 - PaxDollar
$ws.Range("E37").Value = "  -0.56%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.95%  "

# Row 40 - ARBITRUM
$ws.Range("D40").Value = "'0.788"
$ws.Range("E40").Value = "  +1.13%  "

# Row 41 - ImmutableX
$ws.Range("D41").Value = "'0.499"
$ws.Range("E41").Value = "  +1.06%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "'0.774"
$ws.Range("E42").Value = "  +2.94%  "

# Row 43 - RocketPoolETH
$ws.Range("D43").Value = "1.741.97"
$ws.Range("E43").Value = "  +2.17%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  +0.78%  "

# Row 45 - Quant
$ws.Range("D45").Value = "'92.71"
$ws.Range("E45").Value = "  -0.98%  "

# Row 46 - BabyDogeCoin (subscript digit keeps it non-numeric already)
$ws.Range("D46").Value = "0.0$([char]0x2086)0114"
$ws.Range("E46").Value = "  +2.52%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "'1.55"
$ws.Range("E47").Value = "  +9.24%  "

# Row 48 - Aave
$ws.Range("D48").Value = "'53.51"
$ws.Range("E48").Value = "  +1.31%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  -0.14%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  +0.33%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  -0.45%  "
